$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Capture the existing data rows (2-6) before shifting them down by one row
$a2 = $ws.Cells.Item(2,1).Value2
$b2 = $ws.Cells.Item(2,2).Value2
$a3 = $ws.Cells.Item(3,1).Value2
$b3 = $ws.Cells.Item(3,2).Value2
$a4 = $ws.Cells.Item(4,1).Value2
$b4 = $ws.Cells.Item(4,2).Value2
$a5 = $ws.Cells.Item(5,1).Value2
$b5 = $ws.Cells.Item(5,2).Value2
$a6 = $ws.Cells.Item(6,1).Value2
$b6 = $ws.Cells.Item(6,2).Value2

# Remove the existing hyperlinks so they don't get dragged along with the values
$ws.Hyperlinks.Delete()

# Shift rows 2-6 down into rows 3-7
$ws.Cells.Item(7,1).Value = $a6
$ws.Cells.Item(7,2).Value = $b6
$ws.Cells.Item(6,1).Value = $a5
$ws.Cells.Item(6,2).Value = $b5
$ws.Cells.Item(5,1).Value = $a4
$ws.Cells.Item(5,2).Value = $b4
$ws.Cells.Item(4,1).Value = $a3
$ws.Cells.Item(4,2).Value = $b3
$ws.Cells.Item(3,1).Value = $a2
$ws.Cells.Item(3,2).Value = $b2

# Write the new row 2 data
$ws.Cells.Item(2,1).Value = "2076-3387"
$ws.Cells.Item(2,2).Value = "http://susy.mdpi.com"

# Recreate the hyperlinks in the same relative order as the original commit:
# existing ones first (now one row further down), then the newly-added one
$ws.Hyperlinks.Add($ws.Cells.Item(3,2), "http://ade.sagepub.com/submission")
$ws.Cells.Item(3,2).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(4,2), "http://www.mdpi.com/journal/agronomy/submission")
$ws.Cells.Item(4,2).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(5,2), "http://www.aimspress.com/journal/Materials/submission")
$ws.Cells.Item(5,2).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(7,2), "http://campus.usal.es/~revistas_trabajo/index.php/1130-2887/submission")
$ws.Cells.Item(7,2).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(2,2), "http://susy.mdpi.com")
$ws.Cells.Item(2,2).Style = "Hyperlink"

$ws.Range("B16").Select()
